# Scheduled runner update: refresh market-price-derived columns (H-N)
# on the Leve profit sheets. Values are sourced from the latest
# Universalis snapshot; K/L/M/N are re-derived from the new H/I/J prices.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2000
$ws.Range("I40").Value = 2500
$ws.Range("J40").Value = 1500
$ws.Range("K40").Value = 2500
$ws.Range("L40").Value = 1500
$ws.Range("M40").Value = -2325
$ws.Range("N40").Value = -1850

$ws.Range("H64").Value = 4416
$ws.Range("J64").Value = 4665.5
$ws.Range("L64").Value = 4665.5
$ws.Range("N64").Value = -5161.5

$ws.Range("H67").Value = 4416
$ws.Range("J67").Value = 4665.5
$ws.Range("L67").Value = 4665.5
$ws.Range("N67").Value = -6381.5

$ws.Range("H98").Value = 932.3333
$ws.Range("I98").Value = 758.5789
$ws.Range("J98").Value = 2583
$ws.Range("K98").Value = 758.5789
$ws.Range("L98").Value = 2583
$ws.Range("M98").Value = 739.4211
$ws.Range("N98").Value = -5579

$ws.Range("H122").Value = 932.3333
$ws.Range("I122").Value = 758.5789
$ws.Range("J122").Value = 2583
$ws.Range("K122").Value = 2275.7367
$ws.Range("L122").Value = 7749
$ws.Range("M122").Value = 174.2633000000001
$ws.Range("N122").Value = -12649

$ws.Range("H132").Value = 1682.7778
$ws.Range("I132").Value = 1680.625
$ws.Range("K132").Value = 5041.875
$ws.Range("M132").Value = -2511.875

$ws.Range("H135").Value = 700.8
$ws.Range("I135").Value = 506.25
$ws.Range("K135").Value = 4556.25
$ws.Range("M135").Value = -2021.25


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3512.55
$ws.Range("I32").Value = 3512.55
$ws.Range("K32").Value = 3512.55
$ws.Range("M32").Value = -3225.55

$ws.Range("H33").Value = 29876.75
$ws.Range("I33").Value = 19198
$ws.Range("J33").Value = 47674.668
$ws.Range("K33").Value = 19198
$ws.Range("L33").Value = 47674.668
$ws.Range("M33").Value = -18869
$ws.Range("N33").Value = -48332.668

$ws.Range("H45").Value = 2068.6667
$ws.Range("I45").Value = 1982.4
$ws.Range("K45").Value = 1982.4
$ws.Range("M45").Value = -1605.4


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 483.4
$ws.Range("I20").Value = 500.25
$ws.Range("K20").Value = 500.25
$ws.Range("M20").Value = -253.25

$ws.Range("H99").Value = 1317
$ws.Range("I99").Value = 1347.2307
$ws.Range("J99").Value = 1120.5
$ws.Range("K99").Value = 1347.2307
$ws.Range("L99").Value = 1120.5
$ws.Range("M99").Value = 150.7692999999999
$ws.Range("N99").Value = -4116.5


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 501.0345
$ws.Range("I19").Value = 358.15384
$ws.Range("J19").Value = 1739.3334
$ws.Range("K19").Value = 358.15384
$ws.Range("L19").Value = 1739.3334
$ws.Range("M19").Value = -188.15384
$ws.Range("N19").Value = -2079.3334

$ws.Range("H24").Value = 501.0345
$ws.Range("I24").Value = 358.15384
$ws.Range("J24").Value = 1739.3334
$ws.Range("K24").Value = 358.15384
$ws.Range("L24").Value = 1739.3334
$ws.Range("M24").Value = -188.15384
$ws.Range("N24").Value = -2079.3334

$ws.Range("H25").Value = 1527.75
$ws.Range("I25").Value = 1527.75
$ws.Range("K25").Value = 1527.75
$ws.Range("M25").Value = -1353.75

$ws.Range("H31").Value = 2837.8708
$ws.Range("I31").Value = 2021.2593
$ws.Range("K31").Value = 2021.2593
$ws.Range("M31").Value = -1726.2593

$ws.Range("H34").Value = 2837.8708
$ws.Range("I34").Value = 2021.2593
$ws.Range("K34").Value = 2021.2593
$ws.Range("M34").Value = -1819.2593

$ws.Range("H58").Value = 2942
$ws.Range("I58").Value = 3122.6667
$ws.Range("K58").Value = 3122.6667
$ws.Range("M58").Value = -2919.6667

$ws.Range("H92").Value = 35516
$ws.Range("J92").Value = 36719.2
$ws.Range("L92").Value = 36719.2
$ws.Range("N92").Value = -41711.2

$ws.Range("H99").Value = 7999
$ws.Range("I99").Value = 7999
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 7999
$ws.Range("L99").Value = 0
$ws.Range("M99").ClearContents()
$ws.Range("N99").Value = -6501

$ws.Range("H125").Value = 49999.5
$ws.Range("J125").Value = 49999.5
$ws.Range("L125").Value = 49999.5
$ws.Range("N125").Value = -54919.5

$ws.Range("H126").Value = 7999
$ws.Range("I126").Value = 7999
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 23997
$ws.Range("L126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("N126").Value = -21527

$ws.Range("H132").Value = 980
$ws.Range("I132").Value = 980
$ws.Range("K132").Value = 2940
$ws.Range("M132").Value = -410

$ws.Range("H134").Value = 2412.7
$ws.Range("I134").Value = 2275.7896
$ws.Range("K134").Value = 6827.3688
$ws.Range("M134").Value = -4292.3688

$ws.Range("H136").Value = 2942
$ws.Range("I136").Value = 3122.6667
$ws.Range("K136").Value = 9368.000100000001
$ws.Range("M136").Value = -6818.000100000001


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 549
$ws.Range("J38").Value = 1497.5
$ws.Range("L38").Value = 4492.5
$ws.Range("N38").Value = -5186.5

$ws.Range("H60").Value = 758.4
$ws.Range("I60").Value = 464
$ws.Range("K60").Value = 1392
$ws.Range("M60").Value = -1141

$ws.Range("H117").Value = 954.75
$ws.Range("I117").Value = 809.5
$ws.Range("K117").Value = 2428.5
$ws.Range("M117").Value = 1013.5

$ws.Range("H129").Value = 995.4
$ws.Range("J129").Value = 1132.3334
$ws.Range("L129").Value = 3397.0002
$ws.Range("N129").Value = -13397.0002

$ws.Range("H140").Value = 1686.6923
$ws.Range("I140").Value = 1368.9166
$ws.Range("K140").Value = 4106.7498
$ws.Range("M140").Value = 1073.2502

$ws.Range("H141").Value = 1443
$ws.Range("I141").Value = 1443
$ws.Range("K141").Value = 4329
$ws.Range("M141").Value = 851


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 325.44446
$ws.Range("I2").Value = 366.5
$ws.Range("K2").Value = 366.5
$ws.Range("M2").Value = -253.5

$ws.Range("H14").Value = 13003.5
$ws.Range("J14").Value = 16005
$ws.Range("L14").Value = 16005
$ws.Range("N14").Value = -16341


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3987
$ws.Range("I7").Value = 3987
$ws.Range("K7").Value = 3987
$ws.Range("M7").Value = -3875

$ws.Range("H46").Value = 3569.5151
$ws.Range("I46").Value = 2707.0715
$ws.Range("K46").Value = 2707.0715
$ws.Range("M46").Value = -2519.0715

$ws.Range("H93").Value = 1350
$ws.Range("I93").Value = 1200
$ws.Range("K93").Value = 1200
$ws.Range("M93").Value = 48

$ws.Range("H126").Value = 3987
$ws.Range("I126").Value = 3987
$ws.Range("K126").Value = 11961
$ws.Range("M126").Value = -9491

$ws.Range("H132").Value = 6191.3076
$ws.Range("I132").Value = 5747.8335
$ws.Range("K132").Value = 17243.5005
$ws.Range("M132").Value = -14713.5005

$ws.Range("H136").Value = 3086
$ws.Range("I136").Value = 3003.2
$ws.Range("J136").Value = 3500
$ws.Range("K136").Value = 9009.599999999999
$ws.Range("L136").Value = 10500
$ws.Range("M136").Value = -6459.599999999999
$ws.Range("N136").Value = -15600


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H7").Value = 418.33334
$ws.Range("I7").Value = 1000
$ws.Range("J7").Value = 127.5
$ws.Range("K7").Value = 1000
$ws.Range("L7").Value = 127.5
$ws.Range("M7").Value = -887
$ws.Range("N7").Value = -353.5

$ws.Range("H68").Value = 35000
$ws.Range("J68").Value = 35000
$ws.Range("L68").Value = 35000
$ws.Range("N68").Value = -36622

$ws.Range("H69").Value = 23254.2
$ws.Range("J69").Value = 23254.2
$ws.Range("L69").Value = 23254.2
$ws.Range("N69").Value = -24752.2

$ws.Range("H71").Value = 35000
$ws.Range("J71").Value = 35000
$ws.Range("L71").Value = 105000
$ws.Range("N71").Value = -113112

$ws.Range("H72").Value = 23254.2
$ws.Range("J72").Value = 23254.2
$ws.Range("L72").Value = 69762.60000000001
$ws.Range("N72").Value = -77250.60000000001

$ws.Range("H126").Value = 1066.5294
$ws.Range("I126").Value = 945.6875
$ws.Range("K126").Value = 2837.0625
$ws.Range("M126").Value = -367.0625

$ws.Range("H136").Value = 10277.714
$ws.Range("I136").Value = 9593
$ws.Range("K136").Value = 28779
$ws.Range("M136").Value = -26229
